$d = $word.ActiveDocument

# Replace "double tapped" with "just tapped"
$d.Content.Find.Execute("double tapped", $true, $false, $false, $false, $false,
                         $true, 1, $false, "just tapped", 2)

# Add a comma after "button down/up" and before " it's a click )"
$d.Content.Find.Execute("button down/up it" + [char]8217 + "s a click", $true, $false, $false, $false, $false,
                         $true, 1, $false, "button down/up, it" + [char]8217 + "s a click", 2)
